$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Designator"
$ws.Cells.Item(1, 2).Value = "Mid X"
$ws.Cells.Item(1, 3).Value = "Mid Y"
$ws.Cells.Item(1, 4).Value = "Rotation"
$ws.Cells.Item(1, 5).Value = "Layer"

# Data rows 2-122: Row|Designator|MidX|MidY|Rotation|Layer
$data = @"
2|C1|108.7565|-130.84169299999999|0|top
3|C2|113.8965|-123.24369299999999|0|top
4|C3|113.84650000000001|-130.14369300000001|180|top
5|C4|108.7565|-128.56839299999999|180|top
6|C5|106.2165|-121.18969300000001|-90|top
7|C6|108.7565|-126.269693|180|top
8|C7|103.93049999999999|-121.18969300000001|90|top
9|C8|108.7565|-133.12769299999999|180|top
10|C9|108.7565|-123.983693|0|top
11|C10|86.846500000000006|-119.343693|-90|top
12|C11|108.7565|-135.41369299999999|180|top
13|C12|96.8185|-122.205693|90|top
14|C13|77.396500000000003|-125.043693|180|top
15|C14|77.296499999999995|-119.24369299999999|0|top
16|C15|147.79650000000001|-138.74369300000001|-90|top
17|C16|65.6965|-140.79369299999999|90|top
18|C17|150.59649999999999|-138.54369299999999|90|top
19|C18|148.29650000000001|-143.24369300000001|90|top
20|C19|166.69|-131.93000000000001|90|top
21|C20|91.484499999999997|-126.396693|0|top
22|C21|123.3965|-105.343693|-90|top
23|C22|69.6965|-140.093693|-90|top
24|C23|69.6965|-143.89369300000001|90|top
25|C24|101.64449999999999|-121.18969300000001|90|top
26|C25|91.484499999999997|-134.016693|0|top
27|C26|91.738500000000002|-137.699693|90|top
28|C27|82.594499999999996|-127.793693|0|top
29|C28|102.02549999999999|-137.699693|90|top
30|C29|91.484499999999997|-131.60369299999999|0|top
31|C30|91.484499999999997|-128.80969300000001|180|top
32|C31|104.1845|-137.699693|90|top
33|C32|99.358500000000006|-142.652693|90|top
34|C33|82.594499999999996|-134.14369300000001|180|top
35|C34|194.95500000000001|-129.27500000000001|90|top
36|C35|171.012|-124.52|0|top
37|C36|174.36500000000001|-116.48999999999999|0|top
38|C37|194.97999999999999|-117.88|0|top
39|C38|187.22999999999999|-114.33|0|top
40|C39|124.4965|-113.74369299999999|-90|top
41|C40|108.9965|-113.943693|0|top
42|C41|116.8965|-105.343693|-90|top
43|C42|73.496499999999997|-115.843693|0|top
44|C43|170.63|-98.340000000000003|0|top
45|C44|171|-88.890000000000001|0|top
46|C45|157.19999999999999|-98.340000000000003|0|top
47|C46|157.75999999999999|-88.890000000000001|0|top
48|D1|99.358500000000006|-120.681693|-90|top
49|D2|88.436499999999995|-126.52369299999999|90|top
50|D3|85.896500000000003|-126.52369299999999|-90|top
51|D4|99.485500000000002|-138.33469299999999|90|top
52|D5|88.436499999999995|-135.66769300000001|90|top
53|D6|85.896500000000003|-135.66769300000001|-90|top
54|D7|63.596499999999999|-125.943693|90|top
55|D8|72.141499999999994|-143.94319300000001|-90|top
56|D9|74.554500000000004|-143.94319300000001|-90|top
57|D10|70.796499999999995|-124.543693|90|top
58|D11|170.78999999999999|-93.370000000000005|-90|top
59|D12|157.28|-93.392499999999998|-90|top
60|J1|64.236500000000007|-113.643693|-90|top
61|J4|146.64993699999999|-68.622200000000007|180|top
62|J6|173.47|-129.02000000000001|90|top
63|J9|151.159559|-61.5|180|top
64|J10|190.08000000000001|-52.159999999999997|90|top
65|J11|67.519999999999996|-88.694999999999993|-90|top
66|J12|60.223999999999997|-67.349999999999994|-90|top
67|L1|92.754499999999993|-120.427693|90|top
68|L2|95.802499999999995|-139.09669299999999|90|top
69|L3|175.49000000000001|-98.340000000000003|0|top
70|L4|162|-98.340000000000003|0|top
71|Q1|67.046499999999995|-125.793693|0|top
72|Q3|86.896500000000003|-105.568693|90|top
73|Q4|86.796499999999995|-112.543693|-90|top
74|R1|85.996499999999997|-109.068693|-90|top
75|R2|72.6965|-119.643693|-90|top
76|R3|72.246499999999997|-140.343693|-90|top
77|R4|108.9965|-111.943693|180|top
78|R5|94.496499999999997|-106.74369299999999|180|top
79|R6|108.9965|-110.143693|180|top
80|R7|72.6965|-108.343693|90|top
81|R8|74.554500000000004|-140.343693|-90|top
82|R9|71.4465|-128.49369300000001|180|top
83|R10|145.79650000000001|-138.943693|-90|top
84|R11|87.596500000000006|-109.068693|90|top
85|R12|146.90000000000001|-75.650000000000006|0|top
86|R13|178.52000000000001|-91.599999999999994|90|top
87|R14|94.496499999999997|-109.943693|0|top
88|R15|108.9965|-108.443693|0|top
89|R16|94.496499999999997|-108.343693|0|top
90|R17|164.66999999999999|-91.599999999999994|90|top
91|R18|107.1777|-139.190293|90|top
92|R19|87.9285|-129.825693|180|top
93|R20|84.626499999999993|-129.825693|180|top
94|R21|67.717500000000001|-140.569693|-90|top
95|R22|67.717500000000001|-143.843693|-90|top
96|R23|87.9285|-132.111693|180|top
97|R24|84.626499999999993|-132.111693|180|top
98|R25|80.920000000000002|-90.719999999999999|180|top
99|R26|196.95500000000001|-129.27500000000001|90|top
100|R27|171.19|-122.59|0|top
101|R28|173.66499999999999|-118.69|0|top
102|R29|195.22999999999999|-115.848|0|top
103|R30|187.45500000000001|-112.43000000000001|0|top
104|R31|65.346500000000006|-122.24369299999999|180|top
105|R32|116.59650000000001|-111.818693|-90|top
106|R33|108.9965|-106.693693|180|top
107|R34|80.909999999999997|-83.030000000000001|0|top
108|RN1|80.831500000000005|-86.570188999999999|0|top
109|TH1|109.2097|-139.190293|-90|top
110|U1|113.8965|-126.74369299999999|-90|top
111|U2|101.7465|-109.30619299999999|-90|top
112|U3|100.12050000000001|-129.57169300000001|180|top
113|U4|73.496499999999997|-112.343693|0|top
114|U5|77.296499999999995|-122.143693|180|top
115|U6|174.66|-91.75|90|top
116|U7|120.2465|-111.943693|-90|top
117|U8|131.82149999999999|-131.79369299999999|180|top
118|U9|79.996499999999997|-110.543693|0|top
119|U10|184.94999999999999|-120.70999999999999|0|top
120|U11|70.496499999999997|-133.48469299999999|-90|top
121|U12|161.16999999999999|-91.75|90|top
122|Y1|120.1465|-105.893693|0|top
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|"
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = $parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [double]$parts[4]
    $ws.Cells.Item($r, 5).Value = $parts[5]
}

$ws.Range("A1:E1").Select()
